$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): add new columns G..J, keep A..F as-is ---
$ws.Cells.Item(1,3).Value  = "oop_IYC_model"
$ws.Cells.Item(1,4).Value  = "premium_HD"
$ws.Cells.Item(1,5).Value  = "IYC_Total_Cost"
$ws.Cells.Item(1,6).Value  = "oop_HDHP_model"
$ws.Cells.Item(1,7).Value  = "premium_HDHP"
$ws.Cells.Item(1,8).Value  = "HD_Total_Cost"
$ws.Cells.Item(1,9).Value  = "HSA"
$ws.Cells.Item(1,10).Value = "total_savings"

# --- Data rows (row 2..11) ---
# Columns: A tier, B risk, C oop_IYC_model, D premium_HD, E IYC_Total_Cost,
#          F oop_HDHP_model, G premium_HDHP, H HD_Total_Cost, I HSA, J total_savings
$data = @(
    @(413.48727416992188, 1584, 1997.4873046875,    667.7208251953125,  588, 1255.7208251953125, 852, 1593.7664794921875),
    @(457.80319213867188, 1584, 2041.80322265625,   810.4893798828125,  588, 1398.4893798828125, 852, 1495.3138427734375),
    @(509.64337158203125, 1584, 2093.643310546875,  970.97662353515625, 588, 1558.9765625,        852, 1386.666748046875),
    @(600.16644287109375, 1584, 2184.16650390625,   1246.7025146484375, 588, 1834.7025146484375,  852, 1201.4639892578125),
    @(970.81707763671875, 1584, 2554.81689453125,   2179.317626953125,  588, 2767.317626953125,   852, 639.499267578125),
    @(1037.4964599609375, 3948, 4985.49609375,      2470.68701171875,  1464, 3934.68701171875,   1704, 2754.80908203125),
    @(1171.10205078125,   3948, 5119.10205078125,   2799.93701171875,  1464, 4263.93701171875,   1704, 2559.1650390625),
    @(1285.0806884765625, 3948, 5233.08056640625,   3076.21533203125,  1464, 4540.21533203125,   1704, 2396.865234375),
    @(1460.103759765625,  3948, 5408.103515625,     3487.163818359375, 1464, 4951.1640625,       1704, 2160.939453125),
    @(2089.96630859375,   3948, 6037.96630859375,   4724.49169921875,  1464, 6188.49169921875,   1704, 1553.474609375)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $fmt = $ws.Cells.Item($r,2).NumberFormat

    for ($c = 3; $c -le 10; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $cell.Value = $row[$c - 3]
        $cell.NumberFormat = $fmt
    }
}

$ws.Range("A1:J11").Columns.AutoFit() | Out-Null
